$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text-type cells (some values look numeric, e.g. "1.003" or
# "0.3710") so Excel does not silently coerce them to Number and drop
# formatting such as trailing zeros. A leading apostrophe forces text
# entry without the apostrophe becoming part of the stored value.

$ws.Range("D2").Value = "'" + '26.992.70'
$ws.Range("E2").Value = "'" + '  -0.24%  '

$ws.Range("D3").Value = "'" + '1.826.19'
$ws.Range("E3").Value = "'" + '  +0.23%  '

$ws.Range("D4").Value = "'" + '1.004'
$ws.Range("E4").Value = "'" + '  -0.85%  '

$ws.Range("D5").Value = "'" + '311.68'
$ws.Range("E5").Value = "'" + '  +0.03%  '

$ws.Range("D6").Value = "'" + '1.003'
$ws.Range("E6").Value = "'" + '  -0.78%  '

$ws.Range("D7").Value = "'" + '0.4629'
$ws.Range("E7").Value = "'" + '  -0.21%  '

$ws.Range("D8").Value = "'" + '0.3710'
$ws.Range("E8").Value = "'" + '  +2.31%  '

$ws.Range("D9").Value = "'" + '0.07344'
$ws.Range("E9").Value = "'" + '  +0.78%  '

$ws.Range("D10").Value = "'" + '0.8773'
$ws.Range("E10").Value = "'" + '  +1.31%  '

$ws.Range("D11").Value = "'" + '0.07897'
$ws.Range("E11").Value = "'" + '  +3.79%  '

$ws.Range("D12").Value = "'" + '19.75'
$ws.Range("E12").Value = "'" + '  -0.48%  '

$ws.Range("D13").Value = "'" + '1.820.66'
$ws.Range("E13").Value = "'" + '  -0.55%  '

$ws.Range("D14").Value = "'" + '5.339'
$ws.Range("E14").Value = "'" + '  +0.15%  '

$ws.Range("D15").Value = "'" + '6.553'
$ws.Range("E15").Value = "'" + '  +1.49%  '

$ws.Range("D16").Value = "'" + '91.29'
$ws.Range("E16").Value = "'" + '  -1.52%  '

$ws.Range("E17").Value = "'" + '  -0.61%  '

$ws.Range("D18").Value = "'" + '0.000008845'
$ws.Range("E18").Value = "'" + '  +2.62%  '

$ws.Range("E19").Value = "'" + '  -0.71%  '

$ws.Range("D20").Value = "'" + '14.82'
$ws.Range("E20").Value = "'" + '  +2.49%  '

$ws.Range("D21").Value = "'" + '27.016.42'
$ws.Range("E21").Value = "'" + '  -1.30%  '

$ws.Range("D22").Value = "'" + '5.097'
$ws.Range("E22").Value = "'" + '  -1.61%  '

$ws.Range("D23").Value = "'" + '10.53'
$ws.Range("E23").Value = "'" + '  -0.16%  '

$ws.Range("D24").Value = "'" + '2.071.80'
$ws.Range("E24").Value = "'" + '  -0.87%  '

$ws.Range("D25").Value = "'" + '153.00'
$ws.Range("E25").Value = "'" + '  +1.05%  '

$ws.Range("E26").Value = "'" + '  -0.92%  '

$ws.Range("D27").Value = "'" + '18.40'
$ws.Range("E27").Value = "'" + '  +1.02%  '

$ws.Range("D28").Value = "'" + '2.044'
$ws.Range("E28").Value = "'" + '  -2.58%  '

$ws.Range("D29").Value = "'" + '5.130'
$ws.Range("E29").Value = "'" + '  +1.19%  '

$ws.Range("D30").Value = "'" + '115.76'
$ws.Range("E30").Value = "'" + '  -0.13%  '

$ws.Range("D31").Value = "'" + '0.08877'
$ws.Range("E31").Value = "'" + '  -0.23%  '

$ws.Range("D32").Value = "'" + '2.955'
$ws.Range("E32").Value = "'" + '  -0.26%  '

$ws.Range("D33").Value = "'" + '0.7288'
$ws.Range("E33").Value = "'" + '  -0.03%  '

$ws.Range("D34").Value = "'" + '4.440'
$ws.Range("E34").Value = "'" + '  +0.12%  '

$ws.Range("D35").Value = "'" + '1.132'
$ws.Range("E35").Value = "'" + '  -0.44%  '

$ws.Range("D36").Value = "'" + '2.479'
$ws.Range("E36").Value = "'" + '  -2.22%  '

$ws.Range("D37").Value = "'" + '0.01947'
$ws.Range("E37").Value = "'" + '  +1.88%  '

$ws.Range("D38").Value = "'" + '1.067'
$ws.Range("E38").Value = "'" + '  -0.38%  '

$ws.Range("D39").Value = "'" + '0.05226'
$ws.Range("E39").Value = "'" + '  -0.44%  '

$ws.Range("D40").Value = "'" + '2.949'
$ws.Range("E40").Value = "'" + '  +0.16%  '

$ws.Range("D41").Value = "'" + '7.096'
$ws.Range("E41").Value = "'" + '  -0.20%  '

$ws.Range("D42").Value = "'" + '0.5171'
$ws.Range("E42").Value = "'" + '  -0.53%  '

$ws.Range("D43").Value = "'" + '0.1624'
$ws.Range("E43").Value = "'" + '  -0.31%  '

$ws.Range("B44").Value = "'" + 'Decentraland'
$ws.Range("C44").Value = "'" + 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D44").Value = "'" + '0.4843'
$ws.Range("E44").Value = "'" + '  -0.14%  '

$ws.Range("B45").Value = "'" + 'Aptos'
$ws.Range("C45").Value = "'" + 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D45").Value = "'" + '8.164'
$ws.Range("E45").Value = "'" + '  -0.74%  '

$ws.Range("B46").Value = "'" + 'PaxDollar'
$ws.Range("C46").Value = "'" + 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D46").Value = "'" + '1.003'
$ws.Range("E46").Value = "'" + '  -0.76%  '

$ws.Range("B47").Value = "'" + 'EnergySwap'
$ws.Range("C47").Value = "'" + 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").Value = "'" + '10.19'
$ws.Range("E47").Value = "'" + '  +1.25%  '

$ws.Range("D48").Value = "'" + '102.57'
$ws.Range("E48").Value = "'" + '  -0.60%  '

$ws.Range("D49").Value = "'" + '1.629'
$ws.Range("E49").Value = "'" + '  -0.54%  '

$ws.Range("D50").Value = "'" + '0.06200'
$ws.Range("E50").Value = "'" + '  -0.77%  '

$ws.Range("D51").Value = "'" + '64.78'
$ws.Range("E51").Value = "'" + '  +0.65%  '
